$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = "MECit547"
$ws.Range("B2").Value = 23072604
$ws.Range("C2").Value = "sffmqsl59"
$ws.Range("D2").Value = 'Nc6q!T$4'
$ws.Range("F2").Value = "MajCmAme"
$ws.Range("G2").Value = "yZFB"
